$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(2, 5, 8),
    @(3, 5, 8),
    @(4, 5, 8),
    @(5, 5, 10),
    @(5, 6, 20260104),
    @(6, 5, 8),
    @(7, 5, 10),
    @(7, 6, 20260104),
    @(8, 5, 8),
    @(9, 5, 10),
    @(9, 6, 20260104),
    @(10, 5, 1),
    @(11, 5, 8),
    @(12, 5, 10),
    @(12, 6, 20260104),
    @(13, 5, 8),
    @(14, 5, 8),
    @(15, 5, 8),
    @(16, 5, 4),
    @(17, 5, 10),
    @(17, 6, 20260104),
    @(18, 5, 3),
    @(19, 5, 3),
    @(20, 5, 3),
    @(21, 5, 3),
    @(22, 5, 10),
    @(22, 6, 20260104),
    @(23, 5, 10),
    @(23, 6, 20260104),
    @(24, 5, 10),
    @(24, 6, 20260104),
    @(25, 5, 10),
    @(25, 6, 20260104),
    @(26, 5, 10),
    @(26, 6, 20260104),
    @(27, 5, 2),
    @(28, 5, 3),
    @(29, 5, 3),
    @(30, 5, 3),
    @(31, 5, 3),
    @(32, 5, 3),
    @(33, 5, 3),
    @(34, 5, 3),
    @(35, 5, 3),
    @(37, 5, 3),
    @(38, 5, 3),
    @(39, 5, 3),
    @(40, 5, 1),
    @(41, 5, 1),
    @(42, 5, 3),
    @(43, 5, 10),
    @(43, 6, 20260104),
    @(44, 5, 1),
    @(45, 5, 10),
    @(45, 6, 20260104),
    @(46, 5, 1),
    @(47, 5, 3),
    @(48, 5, 1),
    @(49, 5, 2),
    @(50, 5, 8),
    @(51, 5, 8),
    @(52, 5, 8),
    @(53, 5, 8),
    @(54, 5, 8),
    @(55, 5, 8),
    @(56, 5, 8),
    @(57, 5, 8),
    @(58, 5, 2),
    @(59, 5, 2),
    @(60, 5, 2),
    @(61, 5, 2),
    @(62, 5, 2),
    @(63, 5, 2),
    @(64, 5, 2),
    @(65, 5, 3),
    @(66, 5, 3),
    @(67, 5, 3),
    @(68, 5, 3),
    @(69, 5, 3),
    @(70, 5, 4),
    @(71, 5, 4),
    @(72, 5, 4),
    @(73, 5, 4),
    @(74, 5, 4),
    @(75, 5, 4),
    @(76, 5, 4),
    @(77, 5, 7),
    @(78, 5, 7),
    @(79, 5, 7),
    @(80, 5, 7),
    @(81, 5, 7),
    @(82, 5, 7),
    @(83, 5, 7),
    @(84, 5, 7),
    @(85, 5, 7),
    @(86, 5, 7),
    @(87, 5, 1),
    @(88, 5, 1),
    @(89, 5, 1),
    @(90, 5, 1),
    @(91, 5, 10),
    @(91, 6, 20260104),
    @(92, 5, 1),
    @(93, 5, 7),
    @(94, 5, 4),
    @(95, 5, 6),
    @(96, 5, 4),
    @(97, 5, 4),
    @(98, 5, 4),
    @(99, 5, 4)
)

foreach ($chg in $changes) {
    $ws.Cells.Item($chg[0], $chg[1]).Value = $chg[2]
}
